# Add "Spain" test-data sheet (Zettler Market), based on the existing "Italy" sheet.
$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")
[void]$italy.Activate()

# Duplicate the Italy sheet; the copy is placed right after it.
$italy.Copy($null, $italy)
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# The Italy sheet carries a stray fully-formatted row (row 15) whose styling
# stretches across every column (A:XFD) - drop it so Spain only keeps the 16
# real data rows, like the other market sheets.
$spain.Rows.Item(15).Delete()

# Market-specific values for the new sheet.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2034/T2035/T2036"

# Restore Italy's selection/view state (it's no longer the active tab) and set
# the new sheet's selection/active tab.
[void]$italy.Range("A1:D17").Select()
[void]$spain.Activate()
[void]$spain.Range("D16").Select()
